$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so numeric-looking
# values like "1.00" or "0.0000163" are stored as text, matching
# the original inline-string cell type. Reset the style afterwards
# so cells keep the workbook default style (no explicit style index).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "57.693.09"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "3.082.61"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "516.74"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "141.91"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.434"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").Value = "7.27"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "0.375"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "3.614.04"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "25.70"
$ws.Range("E14").Value = "  -3.84%  "
$ws.Range("D15").Value = "0.0000163"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "57.755.44"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "3.082.78"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "6.10"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").Value = "13.04"
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("D20").Value = "8.10"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "335.37"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "0.502"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "66.09"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").Value = "0.170"
$ws.Range("E25").Value = "  +3.79%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "0.0₃0920"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").Value = "6.37"
$ws.Range("E28").Value = "  -3.87%  "
$ws.Range("D29").Value = "7.16"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").Value = "1.83"
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("D31").Value = "20.87"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").Value = "154.85"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").Value = "27.48"
$ws.Range("E34").Value = "  +7.82%  "
$ws.Range("D35").Value = "4.49"
$ws.Range("E35").Value = "  -3.33%  "
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "1.29"
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("D38").Value = "0.0677"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").Value = "3.131.38"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("D40").Value = "37.07"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").Value = "3.91"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "0.662"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "2.262.44"
$ws.Range("E44").Value = "  +3.03%  "
$ws.Range("E45").Value = "  +5.99%  "
$ws.Range("D46").Value = "1.39"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("D47").Value = "20.08"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "0.935"
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "5.87"
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("D50").Value = "263.38"
$ws.Range("E50").Value = "  +13.99%  "
$ws.Range("D51").Value = "0.0874"
$ws.Range("E51").Value = "  +1.09%  "

# Restore default style on column D so no stray style index is left
# behind by the NumberFormat change above.
$ws.Range("D2:D51").Style = "Normal"

